$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample1")

# Add standard deviation formula below the existing average (row 14),
# leaving row 15 blank to match the existing layout style.
$ws.Range("B16").Formula = "=STDEV(B9:B13)"

# Match the recorded selection state after entering the new formula.
$ws.Range("B13").Select()
